$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r = 12; A = "z0bug.fiscalpos_at_1"; C = "z0bug.fiscalpos_at"; D = "z0bug.tax_22v";  E = "z0bug.tax_eu-1-AT-v" },
    @{ r = 13; A = "z0bug.fiscalpos_at_2"; C = "z0bug.fiscalpos_at"; D = "z0bug.tax_10v";  E = "z0bug.tax_eu-3-AT-v" },
    @{ r = 14; A = "z0bug.fiscalpos_be_1"; C = "z0bug.fiscalpos_be"; D = "z0bug.tax_22v";  E = "z0bug.tax_eu-1-BE-v" },
    @{ r = 15; A = "z0bug.fiscalpos_be_2"; C = "z0bug.fiscalpos_be"; D = "z0bug.tax_10v";  E = "z0bug.tax_eu-3-BE-v" },
    @{ r = 16; A = "z0bug.fiscalpos_nl_1"; C = "z0bug.fiscalpos_nl"; D = "z0bug.tax_22v";  E = "z0bug.tax_eu-1-NL-v" },
    @{ r = 17; A = "z0bug.fiscalpos_nl_2"; C = "z0bug.fiscalpos_nl"; D = "z0bug.tax_10v";  E = "z0bug.tax_eu-3-NL-v" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.A
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
}

$ws.Range("E18").Select()
